$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.024.47'
$ws.Range('E2').Value = '  -2.26%  '
$ws.Range('D3').Value = '2.970.52'
$ws.Range('E3').Value = '  -1.58%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.78'
$ws.Range('E5').Value = '  +1.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.73'
$ws.Range('E6').Value = '  -4.80%  '
$ws.Range('E8').Value = '  -1.72%  '
$ws.Range('D9').Value = '2.969.27'
$ws.Range('E9').Value = '  -1.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.143'
$ws.Range('E10').Value = '  -5.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.75'
$ws.Range('E11').Value = '  +1.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.451'
$ws.Range('E12').Value = '  +1.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000225'
$ws.Range('E13').Value = '  -2.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.89'
$ws.Range('E14').Value = '  -4.43%  '
$ws.Range('E15').Value = '  +1.91%  '
$ws.Range('D16').Value = '3.461.35'
$ws.Range('E16').Value = '  -1.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.98'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').Value = '60.994.90'
$ws.Range('E18').Value = '  -2.26%  '
$ws.Range('D19').Value = '2.968.19'
$ws.Range('E19').Value = '  -1.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '446.37'
$ws.Range('E20').Value = '  -5.55%  '
$ws.Range('E21').Value = '  -1.42%  '
$ws.Range('E22').Value = '  -1.84%  '
$ws.Range('E23').Value = '  -1.34%  '
$ws.Range('E24').Value = '  +0.24%  '
$ws.Range('E25').Value = '  -3.22%  '
$ws.Range('E26').Value = '  -8.73%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.93'
$ws.Range('E28').Value = '  -4.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.64'
$ws.Range('E30').Value = '  +1.01%  '
$ws.Range('E31').Value = '  -5.00%  '
$ws.Range('E32').Value = '  -5.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.07'
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('E34').Value = '  -2.76%  '
$ws.Range('B35').Value = 'PEPE'
$ws.Range('C35').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D35').Value = '0.0₃0781'
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('B36').Value = 'Mantle'
$ws.Range('C36').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.01'
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.69'
$ws.Range('E37').Value = '  -1.91%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.15'
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('B39').Value = 'Cosmos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '9.11'
$ws.Range('E39').Value = '  +1.27%  '
$ws.Range('E40').Value = '  -4.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.119'
$ws.Range('E41').Value = '  +4.73%  '
$ws.Range('E42').Value = '  -9.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '388.38'
$ws.Range('E43').Value = '  -7.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0351'
$ws.Range('E44').Value = '  -1.42%  '
$ws.Range('D45').Value = '2.689.06'
$ws.Range('E45').Value = '  -4.30%  '
$ws.Range('E46').Value = '  -6.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '37.11'
$ws.Range('E47').Value = '  -2.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.67'
$ws.Range('E48').Value = '  +3.29%  '
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('E51').Value = '  -0.22%  '
